$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WEO_Data")

# Capture the values/labels currently sitting in column T (they need to
# move over to column V once the two new year columns are inserted).
$t1Value = $ws.Range("T1").Value2
$t2Value = $ws.Range("T2").Value2

# Move the old T1/T2 content ("Estimates Start After" label and the 2016
# "estimates start after" year) two columns to the right, into V1/V2.
$ws.Range("V1").Value = $t1Value
$ws.Range("V2").Value = $t2Value

# Fill the newly introduced 2020/2021 year columns.
$ws.Range("T1").Value = 2020
$ws.Range("U1").Value = 2021

$ws.Range("T2").Value = 2803.7849999999999
$ws.Range("U2").Value = 2879.777

# Match the numeric formatting style used by the other yearly data cells
# (D2:S2) for the two new data points.
$ws.Range("T2:U2").NumberFormat = $ws.Range("S2").NumberFormat
